# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the Typhon_Profits leve-profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 2635.3333
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2953
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2953
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4201
# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 2635.3333
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2953
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 14765
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21005
# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 917.8
$ws.Range("I106").Value = 917.8
$ws.Range("K106").Value = 917.8
$ws.Range("M106").Value = -286.8
# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 90912730
$ws.Range("I113").Value = 111113550
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 111113550
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -111110296
$ws.Range("N113").Value = -15508
# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 1672
$ws.Range("I127").Value = 745.625
$ws.Range("J127").Value = 2413.1
$ws.Range("K127").Value = 2236.875
$ws.Range("L127").Value = 7239.299999999999
$ws.Range("M127").Value = 2723.125
$ws.Range("N127").Value = -17159.3

$ws = $wb.Worksheets.Item("ARM")
# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 21658.334
$ws.Range("J44").Value = 21658.334
$ws.Range("L44").Value = 21658.334
$ws.Range("N44").Value = -22634.334
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 1699.2667
$ws.Range("I97").Value = 1225.8334
$ws.Range("J97").Value = 3593
$ws.Range("K97").Value = 1225.8334
$ws.Range("L97").Value = 3593
$ws.Range("M97").Value = -729.8334
$ws.Range("N97").Value = -4585
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 16078.086
$ws.Range("I132").Value = 1507.8334
$ws.Range("J132").Value = 103499.6
$ws.Range("K132").Value = 4523.5002
$ws.Range("L132").Value = 310498.8
$ws.Range("M132").Value = -1993.5002
$ws.Range("N132").Value = -315558.8

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2005.3334
$ws.Range("I86").Value = 1798.3334
$ws.Range("J86").Value = 2833.3333
$ws.Range("K86").Value = 1798.3334
$ws.Range("L86").Value = 2833.3333
$ws.Range("M86").Value = -675.3334
$ws.Range("N86").Value = -5079.3333
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2005.3334
$ws.Range("I89").Value = 1798.3334
$ws.Range("J89").Value = 2833.3333
$ws.Range("K89").Value = 8991.666999999999
$ws.Range("L89").Value = 14166.6665
$ws.Range("M89").Value = -3375.666999999999
$ws.Range("N89").Value = -25398.6665
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 2096.75
$ws.Range("I99").Value = 1885.1818
$ws.Range("K99").Value = 1885.1818
$ws.Range("M99").Value = -387.1818000000001
# Row 135 (Leve Item ID 41992)
$ws.Range("H135").Value = 43396
$ws.Range("J135").Value = 43396
$ws.Range("L135").Value = 43396
$ws.Range("N135").Value = -53536
# Row 137 (Leve Item ID 42153)
$ws.Range("H137").Value = 47745
$ws.Range("J137").Value = 47745
$ws.Range("L137").Value = 47745
$ws.Range("N137").Value = -57945

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4084.476
$ws.Range("I31").Value = 4298.375
$ws.Range("J31").Value = 3952.8462
$ws.Range("K31").Value = 4298.375
$ws.Range("L31").Value = 3952.8462
$ws.Range("M31").Value = -4003.375
$ws.Range("N31").Value = -4542.8462
# Row 33 (Leve Item ID 1836)
$ws.Range("H33").Value = 15666.667
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 27000
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 27000
$ws.Range("M33").Value = -9621
$ws.Range("N33").Value = -27758
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4084.476
$ws.Range("I34").Value = 4298.375
$ws.Range("J34").Value = 3952.8462
$ws.Range("K34").Value = 4298.375
$ws.Range("L34").Value = 3952.8462
$ws.Range("M34").Value = -4096.375
$ws.Range("N34").Value = -4356.8462
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 20118916
$ws.Range("I99").Value = 4389029
$ws.Range("J99").Value = 50005700
$ws.Range("K99").Value = 4389029
$ws.Range("L99").Value = 50005700
$ws.Range("M99").Value = -4387531
$ws.Range("N99").Value = -50008696
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 20118916
$ws.Range("I126").Value = 4389029
$ws.Range("J126").Value = 50005700
$ws.Range("K126").Value = 13167087
$ws.Range("L126").Value = 150017100
$ws.Range("M126").Value = -13164617
$ws.Range("N126").Value = -150022040
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 4453.75
$ws.Range("I132").Value = 2979
$ws.Range("K132").Value = 8937
$ws.Range("M132").Value = -6407

$ws = $wb.Worksheets.Item("CUL")
# Row 40 (Leve Item ID 4827)
$ws.Range("H40").Value = 676.6667
$ws.Range("I40").Value = 88.71429000000001
$ws.Range("J40").Value = 1499.8
$ws.Range("K40").Value = 354.85716
$ws.Range("L40").Value = 5999.2
$ws.Range("M40").Value = -285.85716
$ws.Range("N40").Value = -6137.2
# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 1154.5
$ws.Range("I98").Value = 1389.3334
$ws.Range("J98").Value = 450
$ws.Range("K98").Value = 4168.0002
$ws.Range("L98").Value = 1350
$ws.Range("M98").Value = -2670.0002
$ws.Range("N98").Value = -4346
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 704.1667
$ws.Range("J122").Value = 1233.5
$ws.Range("L122").Value = 11101.5
$ws.Range("N122").Value = -16001.5

$ws = $wb.Worksheets.Item("GSM")
# Row 52 (Leve Item ID 4147)
$ws.Range("H52").Value = 16672500
$ws.Range("J52").Value = 16672500
$ws.Range("L52").Value = 16672500
$ws.Range("N52").Value = -16673018
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3889.2
$ws.Range("I80").Value = 3366.8572
$ws.Range("K80").Value = 3366.8572
$ws.Range("M80").Value = -2368.8572
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3889.2
$ws.Range("I83").Value = 3366.8572
$ws.Range("K83").Value = 16834.286
$ws.Range("M83").Value = -11842.286
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 30008.79
$ws.Range("I132").Value = 3583.6428
$ws.Range("J132").Value = 103999.2
$ws.Range("K132").Value = 10750.9284
$ws.Range("L132").Value = 311997.6
$ws.Range("M132").Value = -8220.928400000001
$ws.Range("N132").Value = -317057.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 7000.5
$ws.Range("I22").Value = 7000.5
$ws.Range("K22").Value = 7000.5
$ws.Range("M22").Value = -6705.5
# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 7000.5
$ws.Range("I27").Value = 7000.5
$ws.Range("K27").Value = 7000.5
$ws.Range("M27").Value = -6893.5
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2718.0908
$ws.Range("I46").Value = 2433
$ws.Range("J46").Value = 2825
$ws.Range("K46").Value = 2433
$ws.Range("L46").Value = 2825
$ws.Range("M46").Value = -2245
$ws.Range("N46").Value = -3201
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 240.83333
$ws.Range("I55").Value = 192.66667
$ws.Range("K55").Value = 192.66667
$ws.Range("M55").Value = -19.66667000000001
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 5156.4165
$ws.Range("I68").Value = 2520
$ws.Range("K68").Value = 2520
$ws.Range("M68").Value = -1771
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 5156.4165
$ws.Range("I71").Value = 2520
$ws.Range("K71").Value = 12600
$ws.Range("M71").Value = -8856
# Row 110 (Leve Item ID 25809)
$ws.Range("H110").Value = 3360000
$ws.Range("J110").Value = 3360000
$ws.Range("L110").Value = 3360000
$ws.Range("N110").Value = -3368180

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (Leve Item ID 19981)
$ws.Range("H100").Value = 947.9286
$ws.Range("I100").Value = 627.1
$ws.Range("J100").Value = 1750
$ws.Range("K100").Value = 1254.2
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -713.2
$ws.Range("N100").Value = -4582
# Row 103 (Leve Item ID 18548)
$ws.Range("H103").Value = 19602
$ws.Range("J103").Value = 19602
$ws.Range("L103").Value = 19602
$ws.Range("N103").Value = -21946
# Row 106 (Leve Item ID 19652)
$ws.Range("H106").Value = 26188.5
$ws.Range("J106").Value = 26188.5
$ws.Range("L106").Value = 26188.5
$ws.Range("N106").Value = -28712.5
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 4133612
$ws.Range("I107").Value = 883.8570999999999
$ws.Range("J107").Value = 11365886
$ws.Range("K107").Value = 2651.5713
$ws.Range("L107").Value = 34097658
$ws.Range("M107").Value = -731.5712999999996
$ws.Range("N107").Value = -34101498
# Row 109 (Leve Item ID 27161)
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 3004291
$ws.Range("I113").Value = 1641.4286
$ws.Range("J113").Value = 13513564
$ws.Range("K113").Value = 4924.2858
$ws.Range("L113").Value = 40540692
$ws.Range("M113").Value = -2754.2858
$ws.Range("N113").Value = -40545032
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 19609080
$ws.Range("I136").Value = 26316846
$ws.Range("J136").Value = 1761.8462
$ws.Range("K136").Value = 78950538
$ws.Range("L136").Value = 5285.5386
$ws.Range("M136").Value = -78947988
$ws.Range("N136").Value = -10385.5386

